$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 131064775
$ws.Range("B2").Value = 57884
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 442085
$ws.Range("R2").Value = 7039138
$ws.Range("AC2").Value = "Ringhack"

# Row 3 updates
$ws.Range("A3").Value = 131064773
$ws.Range("Q3").Value = 442108
$ws.Range("AC3").Value = "Ringhack äldre"

# Row 4 updates
$ws.Range("A4").Value = 131064784
$ws.Range("B4").Value = 91829
$ws.Range("E4").Value = 5432
$ws.Range("F4").Value = "Granticka"
$ws.Range("G4").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H4").Value = ""
$ws.Range("Q4").Value = 442100
$ws.Range("R4").Value = 7039221
$ws.Range("AC4").Value = ""

# Row 12
$ws.Range("B12").Value = 91829

# Row 13
$ws.Range("B13").Value = 91805

# Row 17
$ws.Range("B17").Value = 91805

# Row 18
$ws.Range("B18").Value = 91805

# Row 20
$ws.Range("B20").Value = 91805

# Row 21
$ws.Range("B21").Value = 91829
